# Update column F (dSF) values on specific rows to re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    7  = -11
    12 = -2
    14 = 13
    15 = 0
    17 = 3
    30 = -9
    32 = -9
    34 = -6
    35 = -5
    36 = 2
    38 = -3
    39 = -1
    42 = -3
    46 = -4
    48 = -6
    54 = 3
    55 = -2
    56 = -3
    60 = -8
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

$wb.Save()
